$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Week 1 ")

# Update Student Name (col A) / Student Email (col B) for rows 2-11.
# Net effect: assign unique example email addresses per row (same person keeps
# the same email), and normalize a couple of name strings.

$ws.Range("A2").Value = "Sebastian "
$ws.Range("B2").Value = "examplemail1"

$ws.Range("A3").Value = "Jeff "
$ws.Range("B3").Value = "examplemail2"

$ws.Range("A4").Value = "Steve "
$ws.Range("B4").Value = "examplemail3"

$ws.Range("A5").Value = "Kathi "
$ws.Range("B5").Value = "examplemail4"

$ws.Range("A6").Value = "Spandhana "
$ws.Range("B6").Value = "examplemail5"

$ws.Range("A7").Value = "Stefan"
$ws.Range("B7").Value = "examplemail6"

$ws.Range("A8").Value = "Franziska "
$ws.Range("B8").Value = "examplemail7"

$ws.Range("A9").Value = "Sebastian "
$ws.Range("B9").Value = "examplemail1"

$ws.Range("A10").Value = "Stefan "
$ws.Range("B10").Value = "examplemail6"

$ws.Range("A11").Value = "Franziska"
$ws.Range("B11").Value = "examplemail7"

$ws.Range("B11").Select()
